# Project alert1 is saved.TEST Author: admin. Type: SAVE.
# Update rule R30's "From" (min hour) threshold on the Rules sheet: 18 -> 19.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

$ws.Range("C10").Value = 19
